$wb = $excel.ActiveWorkbook

# Fix header on the "EQUIPES" sheet: "EQUIPES" -> "EQUIPE"
$wsEquipes = $wb.Worksheets.Item("EQUIPES")
$wsEquipes.Range("B1").Value = "EQUIPE"

# Fix typo on the "PROJETOS" sheet: "Finaça de mesa" -> "Finança de mesa"
$wsProjetos = $wb.Worksheets.Item("PROJETOS")
$wsProjetos.Range("B4").Value = "Finança de mesa"

# Column B on PROJETOS is best-fit sized; widen it slightly to account for
# the extra character introduced by the corrected spelling.
$wsProjetos.Columns.Item(2).ColumnWidth = 13.8333333
